$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update benchmark values per new GETCONSTARRAY results ---
$ws.Range("H4").Value = 25.1
$ws.Range("L4").Value = 199.6
$ws.Range("O4").Value = 106.2
$ws.Range("L5").Value = 20.399999999999999
$ws.Range("H6").Value = 11.4
$ws.Range("L6").Value = 112.7
$ws.Range("O6").Value = 61.7
$ws.Range("H7").Value = 1.3
$ws.Range("L7").Value = 13.1
$ws.Range("O7").Value = 2.2999999999999998
$ws.Range("H8").Value = 10.4
$ws.Range("L8").Value = 53.2
$ws.Range("O8").Value = 40.1
$ws.Range("H10").Value = 23
$ws.Range("L10").Value = 182.2
$ws.Range("L11").Value = 20.399999999999999
$ws.Range("H12").Value = 8.6999999999999993
$ws.Range("L12").Value = 101
$ws.Range("O12").Value = 49.5
$ws.Range("L13").Value = 11.7
$ws.Range("O13").Value = 3.3
$ws.Range("H14").Value = 10.4
$ws.Range("L14").Value = 49
$ws.Range("O14").Value = 39.799999999999997
$ws.Range("H16").Value = 21
$ws.Range("L16").Value = 168
$ws.Range("O16").Value = 86.7
$ws.Range("L17").Value = 20.399999999999999
$ws.Range("H18").Value = 6
$ws.Range("L18").Value = 85
$ws.Range("O18").Value = 39.700000000000003
$ws.Range("H19").Value = 2.5
$ws.Range("L19").Value = 13.6
$ws.Range("O19").Value = 4.5
$ws.Range("H20").Value = 10.4
$ws.Range("L20").Value = 49
$ws.Range("O20").Value = 39.799999999999997
$ws.Range("H22").Value = 20.100000000000001
$ws.Range("L22").Value = 156.4
$ws.Range("O22").Value = 79.5
$ws.Range("L23").Value = 20.399999999999999
$ws.Range("O23").Value = 2.9
$ws.Range("H24").Value = 5
$ws.Range("L24").Value = 72.099999999999994
$ws.Range("O24").Value = 31
$ws.Range("H25").Value = 2.6
$ws.Range("L25").Value = 14.9
$ws.Range("O25").Value = 5.8
$ws.Range("H26").Value = 10.4
$ws.Range("L26").Value = 49
$ws.Range("O26").Value = 39.799999999999997
$ws.Range("L28").Value = 156.30000000000001
$ws.Range("O28").Value = 74.900000000000006
$ws.Range("L29").Value = 20.399999999999999
$ws.Range("H30").Value = 4.2
$ws.Range("L30").Value = 72.099999999999994
$ws.Range("O30").Value = 25.5
$ws.Range("H31").Value = 2.6
$ws.Range("L31").Value = 14.9
$ws.Range("O31").Value = 6.7
$ws.Range("H32").Value = 10.4
$ws.Range("L32").Value = 49
$ws.Range("O32").Value = 39.799999999999997
$ws.Range("H34").Value = 17.899999999999999
$ws.Range("L34").Value = 156.30000000000001
$ws.Range("O34").Value = 70.8
$ws.Range("L35").Value = 20.399999999999999
$ws.Range("O35").Value = 3.9
$ws.Range("H36").Value = 2.9
$ws.Range("L36").Value = 72.099999999999994
$ws.Range("O36").Value = 20.399999999999999
$ws.Range("H37").Value = 2.6
$ws.Range("L37").Value = 14.9
$ws.Range("O37").Value = 6.8
$ws.Range("H38").Value = 10.4
$ws.Range("L38").Value = 49
$ws.Range("O38").Value = 39.799999999999997
$ws.Range("L40").Value = 156.30000000000001
$ws.Range("O40").Value = 70.5
$ws.Range("L41").Value = 20.399999999999999
$ws.Range("O41").Value = 4.5
$ws.Range("H42").Value = 2.4
$ws.Range("L42").Value = 72
$ws.Range("O42").Value = 19.3
$ws.Range("L43").Value = 14.9
$ws.Range("O43").Value = 6.9
$ws.Range("H44").Value = 10.4
$ws.Range("L44").Value = 49
$ws.Range("O44").Value = 39.799999999999997

# --- Add new "UPDATED 20180301" marker cell with red fill, Q5:R5 ---
$ws.Range("Q5").Value = "UPDATED 20180301"
$ws.Range("Q5:R5").Interior.Color = 255
$ws.Range("Q5:R5").Select()
